$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear all existing cell content/formatting so the shared-string table and
# row layout are rebuilt cleanly in natural (row-by-row) order, then rewrite
# every cell to match the new 18-row layout (rows 19-20 are merged away).
$ws.Cells.Clear()

# Row 1: Total capacity
$ws.Range("A1").Value = "Total capacity"
$ws.Range("B1").Value = 1000.0
$ws.Range("B1").NumberFormat = '#,##0'
$ws.Range("C1").Value = "MWh"

# Row 2: Cost per unit
$ws.Range("A2").Value = "Cost per unit"
$ws.Range("B2").Value = 1000.0
$ws.Range("B2").NumberFormat = '"$"#,##0'
$ws.Range("C2").Value = "/kWh"

# Row 3: Cost of capacity
$ws.Range("A3").Value = "Cost of capacity"
$ws.Range("B3").Formula = "=B2*(B1)"
$ws.Range("B3").NumberFormat = '"$"#,##0'

# Row 4: Lifetime
$ws.Range("A4").Value = "Lifetime"
$ws.Range("B4").Value = 10.0
$ws.Range("B4").NumberFormat = 'General'
$ws.Range("C4").Value = "years"

# Row 5: Rate of return
$ws.Range("A5").Value = "Rate of return"
$ws.Range("B5").Value = 0.08
$ws.Range("B5").NumberFormat = '0%'

# Row 6: Monthly cost
$ws.Range("A6").Value = "Monthly cost"
$ws.Range("B6").Formula = "=-PMT(B5/12,B4*12,B3)"
$ws.Range("B6").NumberFormat = '"$"#,##0.00'
$ws.Range("C6").Value = "/month"

# Row 7: Base case storage price
$ws.Range("A7").Value = "Base case storage price"
$ws.Range("B7").Formula = "=B6/B1"
$ws.Range("B7").NumberFormat = '"$"#,##0.00'
$ws.Range("C7").Value = "/MWh.h"

# Row 8: Reservation storage capacity
$ws.Range("A8").Value = "Reservation storage capacity"
$ws.Range("B8").Value = 0.2
$ws.Range("B8").NumberFormat = '0%'
$ws.Range("C8").Value = "of total"

# Row 9: Study case storage price
$ws.Range("A9").Value = "Study case storage price"
$ws.Range("B9").Formula = "=B7*B8"
$ws.Range("B9").NumberFormat = '"$"#,##0.00'
$ws.Range("C9").Value = "/MWh.h"

# Row 10: Energy Price
$ws.Range("A10").Value = "Energy Price"
$ws.Range("B10").Value = 100.0
$ws.Range("B10").NumberFormat = '"$"#,##0'
$ws.Range("C10").Value = "/MWh"

# Row 11: Storage capacity
$ws.Range("A11").Value = "Storage capacity"
$ws.Range("B11").Formula = "=B12*0.15"
$ws.Range("B11").NumberFormat = 'General'
$ws.Range("C11").Value = "MWh"

# Row 12: Power capacity
$ws.Range("A12").Value = "Power capacity"
$ws.Range("B12").Value = 100.0
$ws.Range("B12").NumberFormat = 'General'
$ws.Range("C12").Value = "MW"

# Row 13: Base case hourly storage cost
$ws.Range("A13").Value = "Base case hourly storage cost"
$ws.Range("B13").Formula = "=B11*B7"
$ws.Range("B13").NumberFormat = '"$"#,##0.00'
$ws.Range("C13").Value = "/h"

# Row 14: Study case hourly storage cost
$ws.Range("A14").Value = "Study case hourly storage cost"
$ws.Range("B14").Formula = "=B11*B9"
$ws.Range("B14").NumberFormat = '"$"#,##0.00'
$ws.Range("C14").Value = "/h"

# Row 15: Hourly energy cost
$ws.Range("A15").Value = "Hourly energy cost"
$ws.Range("B15").Formula = "=B12*B10"
$ws.Range("B15").NumberFormat = '"$"#,##0'
$ws.Range("C15").Value = "/h"

# Row 16: Base case storage price sensitivity (a)
$ws.Range("A16").Value = "Base case storage price sensitivity (a)"
$ws.Range("B16").Formula = "=2*B13/B11^2"
$ws.Range("B16").NumberFormat = '"$"#,##0.00'
$ws.Range("C16").Value = "/MW^2.h^3"

# Row 17: Study case storage price sensitivity (a)
$ws.Range("A17").Value = "Study case storage price sensitivity (a)"
$ws.Range("B17").Formula = "=2*B14/B11^2"
$ws.Range("B17").NumberFormat = '"$"#,##0.00'
$ws.Range("C17").Value = "/MW^2.h^3"

# Row 18: Energy price sensitivity (c)
$ws.Range("A18").Value = "Energy price sensitivity (c)"
$ws.Range("B18").Formula = "=2*B15/B12^2"
$ws.Range("B18").NumberFormat = '"$"#,##0'
$ws.Range("C18").Value = "/MW^2.h"

# Column A width: 12.88 -> 29.5 (ColumnWidth excludes ~0.8333 padding baked into the xml "width")
$ws.Columns.Item(1).ColumnWidth = 28.666666666666668
